{"js": "// cetak susulan dp2nt16 dp3n31 1nt9\n//\n// This label sheet has two identical \"cards\" (it's a 2-up mail-merge\n// layout) and each card holds a handful of MERGEFIELD results cached\n// as plain text runs (NO, NAMA, SEPATU, KAOS, TOPI, UB_1..UB_8). We\n// replace the previous student's cached values with the new one's,\n// leaving GENDER / KELAS untouched (they don't change), and leaving\n// everything else (field codes, other runs, formatting) intact.\n//\n// Because several of the old values are short and/or repeat within a\n// card (e.g. \"M\", \"24\", \"26\"), plain whole-document text search is not\n// safe (it can match inside hidden MERGEFIELD instruction text, or hit\n// the wrong occurrence of a repeated number). So every search below is\n// scoped to the specific paragraph that holds the value, and - where a\n// value is genuinely ambiguous even within that paragraph - we use the\n// unique surrounding label text (\"Kaos: M\") or rely on left-to-right\n// match order within that one paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Each \"card\" starts its info block at a different paragraph index;\n// find them by locating the paragraphs that carry our known anchors\n// instead of hard-coding indices, so the script is resilient to minor\n// structural differences.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\n\nfunction findParaIndexes(predicate) {\n  const idxs = [];\n  for (let i = 0; i < texts.length; i++) {\n    if (predicate(texts[i])) idxs.push(i);\n  }\n  return idxs;\n}\n\nasync function replaceInParagraph(paraIndex, searchText, replacement) {\n  const para = paragraphs.items[paraIndex];\n  const results = para.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Expected to find ${JSON.stringify(searchText)} in paragraph ${paraIndex}`);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nasync function replaceNthInParagraph(paraIndex, searchText, occurrence, replacement) {\n  const para = paragraphs.items[paraIndex];\n  const results = para.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length <= occurrence) {\n    throw new Error(\n      `Expected at least ${occurrence + 1} occurrence(s) of ${JSON.stringify(searchText)} in paragraph ${paraIndex}`\n    );\n  }\n  results.items[occurrence].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Paragraph that still holds the old \"NO\" value (printed as e.g. \"H11\").\nconst noParas = findParaIndexes((t) => t.indexOf(\"H11\") !== -1);\n// Paragraph with the cached \"NAMA\" value.\nconst namaParas = findParaIndexes((t) => t.indexOf(\"DOMINGGUS RADJA\") !== -1);\n// Paragraph with Shoes / Kaos / Topi (SEPATU / KAOS / TOPI).\nconst shoesParas = findParaIndexes((t) => t.indexOf(\"Kaos:\") !== -1);\n// Paragraph with \"Uk. Baju\" (UB_1..UB_8).\nconst ukBajuParas = findParaIndexes((t) => t.indexOf(\"Uk. Baju\") !== -1);\n\nif (noParas.length !== 2 || namaParas.length !== 2 || shoesParas.length !== 2 || ukBajuParas.length !== 2) {\n  throw new Error(\n    `Unexpected card count - NO:${noParas.length} NAMA:${namaParas.length} SHOES:${shoesParas.length} UB:${ukBajuParas.length}`\n  );\n}\n\nfor (const idx of noParas) {\n  await replaceInParagraph(idx, \"H11\", \"H17\");\n}\n\nfor (const idx of namaParas) {\n  await replaceInParagraph(idx, \"DOMINGGUS RADJA\", \"HENDRI RUSMAWARDANA\");\n}\n\nfor (const idx of shoesParas) {\n  await replaceInParagraph(idx, \"40\", \"42\"); // SEPATU\n  await replaceInParagraph(idx, \"Kaos: M\", \"Kaos: XL\"); // KAOS (ambiguous \"M\" on its own)\n  await replaceInParagraph(idx, \"55\", \"58\"); // TOPI\n}\n\n// UB_1 .. UB_8, in left-to-right order: 45;24;16;26;24;26;70;41\n//                                    -> 46;25;18;30;30;30;73;42\nfor (const idx of ukBajuParas) {\n  await replaceInParagraph(idx, \"45\", \"46\"); // UB_1 (unique in paragraph)\n  await replaceNthInParagraph(idx, \"24\", 0, \"25\"); // UB_2 (1st \"24\")\n  await replaceInParagraph(idx, \"16\", \"18\"); // UB_3 (unique in paragraph)\n  await replaceNthInParagraph(idx, \"26\", 0, \"30\"); // UB_4 (1st \"26\")\n  await replaceNthInParagraph(idx, \"24\", 0, \"30\"); // UB_5 (was 2nd \"24\", now 1st remaining)\n  await replaceNthInParagraph(idx, \"26\", 0, \"30\"); // UB_6 (was 2nd \"26\", now 1st remaining)\n  await replaceInParagraph(idx, \"70\", \"73\"); // UB_7 (unique in paragraph)\n  await replaceInParagraph(idx, \"41\", \"42\"); // UB_8 (unique in paragraph)\n}\n", "ps1": "# Cetak susulan dp2nt16 dp3n31 1nt9\n# Update the cached MERGEFIELD results on the \"Setting Baju (Hal depan) F4 PDH\"\n# label sheet (two identical cards on the page) with the new student's data.\n\n$d = $word.ActiveDocument\n\n# Map MERGEFIELD name -> new cached result text.\n# (GENDER and KELAS are intentionally left untouched.)\n$newValues = @{\n    \"NO\"     = \"H17\"\n    \"NAMA\"   = \"HENDRI RUSMAWARDANA\"\n    \"SEPATU\" = \"42\"\n    \"KAOS\"   = \"XL\"\n    \"TOPI\"   = \"58\"\n    \"UB_1\"   = \"46\"\n    \"UB_2\"   = \"25\"\n    \"UB_3\"   = \"18\"\n    \"UB_4\"   = \"30\"\n    \"UB_5\"   = \"30\"\n    \"UB_6\"   = \"30\"\n    \"UB_7\"   = \"73\"\n    \"UB_8\"   = \"42\"\n}\n\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $f = $d.Fields.Item($i)\n    $code = $f.Code.Text.Trim()\n\n    foreach ($name in $newValues.Keys) {\n        if ($code -eq \"MERGEFIELD $name\") {\n            $r = $f.Result\n            $range = $d.Range($r.Start, $r.End)\n            $range.Text = $newValues[$name]\n            break\n        }\n    }\n}\n"}
